# Updated symbol list on Mon Feb 13 05:34:02 UTC 2023 with GitHub Actions
# Refresh Price (column D) and Volume(1h) (column E) values for the
# cryptocurrency rows on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> (new Price, new Volume(1h))
# $null entries mean "leave the existing value unchanged".
$updates = @{
    2  = @("314.22", "2.26%")
    3  = @("40.27", "-1.85%")
    4  = @("5.171", "-1.43%")
    5  = @("0.07580", "-1.14%")
    6  = @("4.318", "0.15%")
    7  = @("1.669", "2.75%")
    8  = @("0.9255", "0.73%")
    10 = @("0.1197", "-4.29%")
    11 = @("0.1820", "-0.34%")
    12 = @("0.09015", "-1.57%")
    13 = @("0.04101", "-3.71%")
    14 = @("0.1053", "0.27%")
    15 = @("0.001278", "1.47%")
    16 = @("0.005800", "1.00%")
    18 = @("3.350", "-0.10%")
    19 = @("0.3356", "0.62%")
    20 = @("7.597", "2.59%")
    21 = @("0.1350", "-2.41%")
    23 = @("0.04027", "-1.04%")
    24 = @("0.001272", "0.81%")
    25 = @("0.004070", "-4.51%")
    26 = @("0.0001271", "-0.16%")
    38 = @($null, "-2.15%")
    39 = @("0.05149", "-2.53%")
    40 = @("0.007722", "-1.72%")
    41 = @($null, "-0.98%")
    42 = @("0.007596", "11.28%")
    43 = @($null, "72.48%")
    44 = @("0.008730", "13.16%")
    45 = @("0.3395", "11.18%")
    46 = @("0.00006587", "-2.05%")
    47 = @($null, "-0.16%")
    48 = @("0.2687", "-11.25%")
    49 = @("0.004200", "2.43%")
    50 = @($null, "-0.16%")
    51 = @($null, "-0.16%")
}

foreach ($row in $updates.Keys) {
    $pair = $updates[$row]
    $price = $pair[0]
    $volume = $pair[1]

    if ($null -ne $price) {
        $cell = $ws.Cells.Item($row, 4)
        $cell.NumberFormat = "@"
        $cell.Value = $price
        $cell.Style = "Normal"
    }
    if ($null -ne $volume) {
        $cell = $ws.Cells.Item($row, 5)
        $cell.NumberFormat = "@"
        $cell.Value = $volume
        $cell.Style = "Normal"
    }
}
